$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks error: Marking row (11) and Total row (12)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "56 / 112"
